# Generate Report for handback
# Update the Correspond Handoff/Handback DateTime values on the
# per-language report sheets (zh-cn, de-de) to reflect a newly
# regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-26 12:30:31"
$wsZhCn.Range("D3").Value = "2016-01-26 12:30:31"
$wsZhCn.Range("G2").Value = "2016-01-26 12:31:21"
$wsZhCn.Range("G3").Value = "2016-01-26 12:31:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-26 12:30:43"
$wsDeDe.Range("D3").Value = "2016-01-26 12:30:43"
$wsDeDe.Range("G2").Value = "2016-01-26 12:31:41"
$wsDeDe.Range("G3").Value = "2016-01-26 12:31:41"
